# PowerBusInfo.xlsx update: v0.0.3r -> v0.0.4r
# Rename the db-keys for the (de)commission year columns from
# comYear/decomYear to YearCom/YearDecom, lowercase the "Excl." db-key
# to "excl", and bump the version string shown at the top of the sheet.
# Apply identically to both scenario sheets (scenarioA / scenarioB).

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    # Bump the version label (row 2, col C)
    $ws.Range("C2").Value2 = "v0.0.4r"

    # Lowercase the "Excl." db-key on the db-key row (row 4)
    $ws.Range("A4").Value2 = "excl"

    # Rename the (de)commission-year db-keys on the db-key row (row 4)
    $ws.Range("K4").Value2 = "YearCom"
    $ws.Range("L4").Value2 = "YearDecom"
}
